$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = 15.4
$ws.Range("L5").Value = 1.5

$ws.Range("F6").Value = 19.4
$ws.Range("L6").Value = 2.1

$ws.Range("F7").Value = 6.4
$ws.Range("L7").Value = 1.2

$ws.Range("F8").Value = 1.5
$ws.Range("L8").Value = 0.8

$ws.Range("F9").Value = 8.300000000000001
$ws.Range("L9").Value = 3.3

$ws.Range("F13").Value = 75.09999999999999
$ws.Range("L13").Value = 3.6

$ws.Range("F14").Value = 19.2
$ws.Range("L14").Value = 6

$ws.Range("F16").Value = 64.60000000000001
$ws.Range("L16").Value = 1.3

$ws.Range("F18").Value = 31.6
$ws.Range("L18").Value = 1.7

$ws.Range("F19").Value = 57.9

$ws.Range("F20").Value = 72.59999999999999
$ws.Range("L20").Value = 0.8999999999999999

$ws.Range("F21").Value = 71.8
$ws.Range("L21").Value = 1.4

$ws.Range("F22").Value = 78.90000000000001

$ws.Range("F23").Value = 93.09999999999999

$ws.Range("F24").Value = 76

$ws.Range("F25").Value = 91.40000000000001

$ws.Range("F26").Value = 38.4

$ws.Range("F27").Value = 0
$ws.Range("L27").Value = 0
